$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 303.46667
$ws.Cells.Item(2, 9).Value = 309.3846
$ws.Cells.Item(2, 10).Value = 265
$ws.Cells.Item(2, 11).Value = 309.3846
$ws.Cells.Item(2, 12).Value = 265
$ws.Cells.Item(2, 13).Value = -196.3846
$ws.Cells.Item(2, 14).Value = -491
# Row 4
$ws.Cells.Item(4, 8).Value = 493
$ws.Cells.Item(4, 9).Value = 400.72726
$ws.Cells.Item(4, 11).Value = 400.72726
$ws.Cells.Item(4, 13).Value = -286.72726
# Row 8
$ws.Cells.Item(8, 8).Value = 361.3
$ws.Cells.Item(8, 9).Value = 361.3
$ws.Cells.Item(8, 11).Value = 1083.9
$ws.Cells.Item(8, 13).Value = -944.9000000000001
# Row 48
$ws.Cells.Item(48, 8).Value = 3746.5
$ws.Cells.Item(48, 10).Value = 3895.8
$ws.Cells.Item(48, 12).Value = 11687.4
$ws.Cells.Item(48, 14).Value = -12271.4
# Row 56
$ws.Cells.Item(56, 8).Value = 3746.5
$ws.Cells.Item(56, 10).Value = 3895.8
$ws.Cells.Item(56, 12).Value = 11687.4
$ws.Cells.Item(56, 14).Value = -12755.4
# Row 70
$ws.Cells.Item(70, 8).Value = 945
$ws.Cells.Item(70, 9).Value = 773.75
$ws.Cells.Item(70, 10).Value = 1042.8572
$ws.Cells.Item(70, 11).Value = 2321.25
$ws.Cells.Item(70, 12).Value = 3128.5716
$ws.Cells.Item(70, 13).Value = -2051.25
$ws.Cells.Item(70, 14).Value = -3668.5716
# Row 73
$ws.Cells.Item(73, 8).Value = 945
$ws.Cells.Item(73, 9).Value = 773.75
$ws.Cells.Item(73, 10).Value = 1042.8572
$ws.Cells.Item(73, 11).Value = 2321.25
$ws.Cells.Item(73, 12).Value = 3128.5716
$ws.Cells.Item(73, 13).Value = -1385.25
$ws.Cells.Item(73, 14).Value = -5000.571599999999
# Row 86
$ws.Cells.Item(86, 8).Value = 4640
$ws.Cells.Item(86, 9).Value = 5050
$ws.Cells.Item(86, 10).Value = 3000
$ws.Cells.Item(86, 11).Value = 5050
$ws.Cells.Item(86, 12).Value = 3000
$ws.Cells.Item(86, 13).Value = -3927
$ws.Cells.Item(86, 14).Value = -5246
# Row 89
$ws.Cells.Item(89, 8).Value = 4640
$ws.Cells.Item(89, 9).Value = 5050
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 25250
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = -19634
$ws.Cells.Item(89, 14).Value = -26232
# Row 96
$ws.Cells.Item(96, 8).Value = 1459.0714
$ws.Cells.Item(96, 9).Value = 1855.3334
$ws.Cells.Item(96, 10).Value = 745.8
$ws.Cells.Item(96, 11).Value = 5566.0002
$ws.Cells.Item(96, 12).Value = 2237.4
$ws.Cells.Item(96, 13).Value = -4193.0002
$ws.Cells.Item(96, 14).Value = -4983.4
# Row 100
$ws.Cells.Item(100, 8).Value = 2031.2858
$ws.Cells.Item(100, 9).Value = 1666.6666
$ws.Cells.Item(100, 10).Value = 2304.75
$ws.Cells.Item(100, 11).Value = 1666.6666
$ws.Cells.Item(100, 12).Value = 2304.75
$ws.Cells.Item(100, 13).Value = -1125.6666
$ws.Cells.Item(100, 14).Value = -3386.75
# Row 106
$ws.Cells.Item(106, 8).Value = 9368.643
$ws.Cells.Item(106, 9).Value = 9993.154
$ws.Cells.Item(106, 11).Value = 9993.154
$ws.Cells.Item(106, 13).Value = -9362.154
# Row 107
$ws.Cells.Item(107, 8).Value = 2801.7144
$ws.Cells.Item(107, 9).Value = 2191.1428
$ws.Cells.Item(107, 10).Value = 4633.4287
$ws.Cells.Item(107, 11).Value = 2191.1428
$ws.Cells.Item(107, 12).Value = 4633.4287
$ws.Cells.Item(107, 13).Value = -271.1428000000001
$ws.Cells.Item(107, 14).Value = -8473.4287
# Row 116
$ws.Cells.Item(116, 8).Value = 2798.1304
$ws.Cells.Item(116, 10).Value = 2688.9092
$ws.Cells.Item(116, 12).Value = 2688.9092
$ws.Cells.Item(116, 14).Value = -9572.9092
# Row 129
$ws.Cells.Item(129, 8).Value = 842.0862
$ws.Cells.Item(129, 10).Value = 863.49054
$ws.Cells.Item(129, 12).Value = 2590.47162
$ws.Cells.Item(129, 14).Value = -12590.47162
# Row 138
$ws.Cells.Item(138, 8).Value = 2159.46
$ws.Cells.Item(138, 9).Value = 1035
$ws.Cells.Item(138, 10).Value = 2423.2222
$ws.Cells.Item(138, 11).Value = 3105
$ws.Cells.Item(138, 12).Value = 7269.6666
$ws.Cells.Item(138, 13).Value = 2035
$ws.Cells.Item(138, 14).Value = -17549.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Cells.Item(132, 8).Value = 2918.875
$ws.Cells.Item(132, 9).Value = 2778.1
$ws.Cells.Item(132, 10).Value = 3341.2
$ws.Cells.Item(132, 11).Value = 8334.299999999999
$ws.Cells.Item(132, 12).Value = 10023.6
$ws.Cells.Item(132, 13).Value = -5804.299999999999
$ws.Cells.Item(132, 14).Value = -15083.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 3725.3713
$ws.Cells.Item(134, 9).Value = 849.65625
$ws.Cells.Item(134, 10).Value = 34399.668
$ws.Cells.Item(134, 11).Value = 2548.96875
$ws.Cells.Item(134, 12).Value = 103199.004
$ws.Cells.Item(134, 13).Value = -13.96875
$ws.Cells.Item(134, 14).Value = -108269.004

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1462.6666
$ws.Cells.Item(31, 9).Value = 1362.8445
$ws.Cells.Item(31, 10).Value = 1961.7778
$ws.Cells.Item(31, 11).Value = 1362.8445
$ws.Cells.Item(31, 12).Value = 1961.7778
$ws.Cells.Item(31, 13).Value = -1067.8445
$ws.Cells.Item(31, 14).Value = -2551.7778
# Row 34
$ws.Cells.Item(34, 8).Value = 1462.6666
$ws.Cells.Item(34, 9).Value = 1362.8445
$ws.Cells.Item(34, 10).Value = 1961.7778
$ws.Cells.Item(34, 11).Value = 1362.8445
$ws.Cells.Item(34, 12).Value = 1961.7778
$ws.Cells.Item(34, 13).Value = -1160.8445
$ws.Cells.Item(34, 14).Value = -2365.7778
# Row 59
$ws.Cells.Item(59, 8).Value = 29500
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 29500
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 29500
$ws.Cells.Item(59, 13).ClearContents()
$ws.Cells.Item(59, 14).Value = -31790
# Row 132
$ws.Cells.Item(132, 8).Value = 2156.8333
$ws.Cells.Item(132, 9).Value = 2260.3845
$ws.Cells.Item(132, 11).Value = 6781.1535
$ws.Cells.Item(132, 13).Value = -4251.1535
# Row 134
$ws.Cells.Item(134, 8).Value = 20001438
$ws.Cells.Item(134, 9).Value = 1512.8182
$ws.Cells.Item(134, 10).Value = 166667550
$ws.Cells.Item(134, 11).Value = 4538.4546
$ws.Cells.Item(134, 12).Value = 500002650
$ws.Cells.Item(134, 13).Value = -2003.4546
$ws.Cells.Item(134, 14).Value = -500007720

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 76
$ws.Cells.Item(76, 8).Value = 5636.364
$ws.Cells.Item(76, 10).Value = 6120
$ws.Cells.Item(76, 12).Value = 18360
$ws.Cells.Item(76, 14).Value = -19126
# Row 79
$ws.Cells.Item(79, 8).Value = 5636.364
$ws.Cells.Item(79, 10).Value = 6120
$ws.Cells.Item(79, 12).Value = 18360
$ws.Cells.Item(79, 14).Value = -21012
# Row 113
$ws.Cells.Item(113, 8).Value = 749.0769
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 749.0769
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 2247.2307
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -6587.2307
# Row 129
$ws.Cells.Item(129, 8).Value = 29762616
$ws.Cells.Item(129, 9).Value = 55556000
$ws.Cells.Item(129, 10).Value = 10417579
$ws.Cells.Item(129, 11).Value = 166668000
$ws.Cells.Item(129, 12).Value = 31252737
$ws.Cells.Item(129, 13).Value = -166663000
$ws.Cells.Item(129, 14).Value = -31262737
# Row 131
$ws.Cells.Item(131, 8).Value = 21309466
$ws.Cells.Item(131, 10).Value = 41572.055
$ws.Cells.Item(131, 12).Value = 124716.165
$ws.Cells.Item(131, 14).Value = -134796.165
# Row 132
$ws.Cells.Item(132, 8).Value = 1259.75
$ws.Cells.Item(132, 9).Value = 1078
$ws.Cells.Item(132, 10).Value = 1320.3334
$ws.Cells.Item(132, 11).Value = 9702
$ws.Cells.Item(132, 12).Value = 11883.0006
$ws.Cells.Item(132, 13).Value = -7172
$ws.Cells.Item(132, 14).Value = -16943.0006
# Row 140
$ws.Cells.Item(140, 8).Value = 32758.152
$ws.Cells.Item(140, 9).Value = 49789.953
$ws.Cells.Item(140, 10).Value = 2952.5
$ws.Cells.Item(140, 11).Value = 149369.859
$ws.Cells.Item(140, 12).Value = 8857.5
$ws.Cells.Item(140, 13).Value = -144189.859
$ws.Cells.Item(140, 14).Value = -19217.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 5392.8335
$ws.Cells.Item(132, 9).Value = 6067.1665
$ws.Cells.Item(132, 11).Value = 18201.4995
$ws.Cells.Item(132, 13).Value = -15671.4995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 5350
$ws.Cells.Item(46, 9).Value = 1500
$ws.Cells.Item(46, 10).Value = 5777.778
$ws.Cells.Item(46, 11).Value = 1500
$ws.Cells.Item(46, 12).Value = 5777.778
$ws.Cells.Item(46, 13).Value = -1312
$ws.Cells.Item(46, 14).Value = -6153.778
# Row 55
$ws.Cells.Item(55, 8).Value = 271.76923
$ws.Cells.Item(55, 9).Value = 217.64516
$ws.Cells.Item(55, 10).Value = 481.5
$ws.Cells.Item(55, 11).Value = 217.64516
$ws.Cells.Item(55, 12).Value = 481.5
$ws.Cells.Item(55, 13).Value = -44.64516
$ws.Cells.Item(55, 14).Value = -827.5
# Row 61
$ws.Cells.Item(61, 8).Value = 1604.125
$ws.Cells.Item(61, 9).Value = 1420.2
$ws.Cells.Item(61, 10).Value = 1910.6666
$ws.Cells.Item(61, 11).Value = 1420.2
$ws.Cells.Item(61, 12).Value = 1910.6666
$ws.Cells.Item(61, 13).Value = -1218.2
$ws.Cells.Item(61, 14).Value = -2314.6666
# Row 113
$ws.Cells.Item(113, 8).Value = 1604.125
$ws.Cells.Item(113, 9).Value = 1420.2
$ws.Cells.Item(113, 10).Value = 1910.6666
$ws.Cells.Item(113, 11).Value = 1420.2
$ws.Cells.Item(113, 12).Value = 1910.6666
$ws.Cells.Item(113, 13).Value = 749.8
$ws.Cells.Item(113, 14).Value = -6250.6666
# Row 132
$ws.Cells.Item(132, 8).Value = 3466.389
$ws.Cells.Item(132, 9).Value = 7332.6665
$ws.Cells.Item(132, 11).Value = 21997.9995
$ws.Cells.Item(132, 13).Value = -19467.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Cells.Item(94, 8).Value = 10950
$ws.Cells.Item(94, 10).Value = 10950
$ws.Cells.Item(94, 12).Value = 10950
$ws.Cells.Item(94, 14).Value = -12752
# Row 107
$ws.Cells.Item(107, 8).Value = 542
$ws.Cells.Item(107, 9).Value = 479
$ws.Cells.Item(107, 11).Value = 1437
$ws.Cells.Item(107, 13).Value = 483
# Row 113
$ws.Cells.Item(113, 8).Value = 731.5
$ws.Cells.Item(113, 9).Value = 398
$ws.Cells.Item(113, 10).Value = 1065
$ws.Cells.Item(113, 11).Value = 1194
$ws.Cells.Item(113, 12).Value = 3195
$ws.Cells.Item(113, 13).Value = 976
$ws.Cells.Item(113, 14).Value = -7535
# Row 132
$ws.Cells.Item(132, 8).Value = 2443.889
$ws.Cells.Item(132, 9).Value = 2370.5356
$ws.Cells.Item(132, 11).Value = 7111.6068
$ws.Cells.Item(132, 13).Value = -4581.6068
